# Add an "Italy" test-data sheet, cloned from the existing "Germany" sheet
# (same layout/styles/merged cells), populated with the Italy market name
# and part number, inserted as the last tab and made the active sheet.
# Also clears the previously-active "Slovakia" tab's selection/active state,
# matching Excel's normal behaviour when the active sheet changes.

$wb = $excel.ActiveWorkbook

$germany  = $wb.Worksheets.Item("Germany")
$slovakia = $wb.Worksheets.Item("Slovakia")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Clone Germany's sheet (formatting, column widths, merged cells, styles)
# and drop the copy after the last existing sheet (i.e. after Slovakia).
$germany.Copy($null, $lastSheet)

$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Italy-specific data
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3145/T2159/NGC-3145/T2159 "

# Slovakia is no longer the active tab; reset its selection to a full-sheet
# selection the way Excel leaves a deactivated sheet.
$slovakia.Activate() | Out-Null
$slovakia.Cells.Select() | Out-Null

# Make the new Italy sheet the active tab, with B4 selected.
$italy.Activate() | Out-Null
$italy.Range("B4").Select() | Out-Null
